$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 27585.55900144784
$ws.Range("C2").Value = 46230.50033102927
$ws.Range("D2").Value = 58753.65503573159
$ws.Range("E2").Value = 60357.44827572616
$ws.Range("F2").Value = 48554.2273382869
$ws.Range("G2").Value = 33639.30172130118
$ws.Range("H2").Value = 23656.21202825921
$ws.Range("B3").Value = 27833.48489556097
$ws.Range("C3").Value = 46537.43871720046
$ws.Range("D3").Value = 59604.10388659112
$ws.Range("E3").Value = 61106.07875624853
$ws.Range("F3").Value = 48749.49792624477
$ws.Range("G3").Value = 34527.19016995178
$ws.Range("H3").Value = 24082.44523951698
$ws.Range("B4").Value = 28120.35798251159
$ws.Range("C4").Value = 46818.99022657244
$ws.Range("D4").Value = 60390.74369802895
$ws.Range("E4").Value = 61687.42132671152
$ws.Range("F4").Value = 48912.35836835435
$ws.Range("G4").Value = 35283.67105294909
$ws.Range("H4").Value = 24402.7160158085
$ws.Range("B5").Value = 28440.1394244515
$ws.Range("C5").Value = 47083.30147352849
$ws.Range("D5").Value = 60982.0133547357
$ws.Range("E5").Value = 62059.56867493501
$ws.Range("F5").Value = 49144.94964935679
$ws.Range("G5").Value = 35969.98615708572
$ws.Range("H5").Value = 24809.75794999196
$ws.Range("B6").Value = 28727.87219284945
$ws.Range("C6").Value = 47365.0056904493
$ws.Range("D6").Value = 61387.64920925194
$ws.Range("E6").Value = 62431.11822533877
$ws.Range("F6").Value = 49517.63539039208
$ws.Range("G6").Value = 36499.6979619663
$ws.Range("H6").Value = 25433.66834417072
$ws.Range("B7").Value = 28765.18674902109
$ws.Range("C7").Value = 47574.45877601231
$ws.Range("D7").Value = 61507.79844878455
$ws.Range("E7").Value = 62800.79723373384
$ws.Range("F7").Value = 50120.52886894477
$ws.Range("G7").Value = 36818.16558707281
$ws.Range("H7").Value = 26180.33003051723
$ws.Range("B8").Value = 28579.26520987642
$ws.Range("C8").Value = 47672.24038638664
$ws.Range("D8").Value = 61216.67099255863
$ws.Range("E8").Value = 63006.55565845839
$ws.Range("F8").Value = 50742.09889504485
$ws.Range("G8").Value = 36989.58790176683
$ws.Range("H8").Value = 26817.46302647276
$ws.Range("B9").Value = 28318.98524040454
$ws.Range("C9").Value = 47609.39139199234
$ws.Range("D9").Value = 60686.76824027306
$ws.Range("E9").Value = 63112.41166332595
$ws.Range("F9").Value = 51296.79272151761
$ws.Range("G9").Value = 37216.92723392993
$ws.Range("H9").Value = 27438.24801797381
$ws.Range("B10").Value = 28095.88653545968
$ws.Range("C10").Value = 47431.41377996354
$ws.Range("D10").Value = 60069.94003182562
$ws.Range("E10").Value = 63266.792995344
$ws.Range("F10").Value = 51611.1552241052
$ws.Range("G10").Value = 37524.67760592084
$ws.Range("H10").Value = 28095.13375800016
$ws.Range("B11").Value = 27935.67792452151
$ws.Range("C11").Value = 47341.88977455212
$ws.Range("D11").Value = 59485.41082347586
$ws.Range("E11").Value = 63652.82040814641
$ws.Range("F11").Value = 51785.63266872131
$ws.Range("G11").Value = 37978.41726048096
$ws.Range("H11").Value = 28633.7908346382
$ws.Range("B12").Value = 27850.84241593184
$ws.Range("C12").Value = 47630.10312957286
$ws.Range("D12").Value = 59096.51125969081
$ws.Range("E12").Value = 64043.21291747566
$ws.Range("F12").Value = 51900.31040329915
$ws.Range("G12").Value = 38579.65515878052
$ws.Range("H12").Value = 29017.63254753576
$ws.Range("B13").Value = 27816.26461778318
$ws.Range("C13").Value = 48263.95805570389
$ws.Range("D13").Value = 59041.99879047045
$ws.Range("E13").Value = 64428.43302098078
$ws.Range("F13").Value = 52116.36595043198
$ws.Range("G13").Value = 39266.75691876117
$ws.Range("H13").Value = 29285.86517269223
$ws.Range("B14").Value = 27924.59699191165
$ws.Range("C14").Value = 49160.55321169037
$ws.Range("D14").Value = 59338.96172045969
$ws.Range("E14").Value = 64841.19283883668
$ws.Range("F14").Value = 52661.53158672046
$ws.Range("G14").Value = 40028.02043507651
$ws.Range("H14").Value = 29628.0149829427
$ws.Range("B15").Value = 28304.19458040293
$ws.Range("C15").Value = 50172.875331264
$ws.Range("D15").Value = 59904.1141107019
$ws.Range("E15").Value = 65268.08011371663
$ws.Range("F15").Value = 53355.42235973399
$ws.Range("G15").Value = 40711.70694506259
$ws.Range("H15").Value = 30068.20548618569
$ws.Range("B16").Value = 29017.19723041427
$ws.Range("C16").Value = 51167.66580301944
$ws.Range("D16").Value = 60705.98464949397
$ws.Range("E16").Value = 65537.25649807356
$ws.Range("F16").Value = 54221.05779432267
$ws.Range("G16").Value = 41382.33378791394
$ws.Range("H16").Value = 30573.73555868528
$ws.Range("B17").Value = 29929.87616018286
$ws.Range("C17").Value = 52175.82728250376
$ws.Range("D17").Value = 61618.46909323441
$ws.Range("E17").Value = 65656.85894519911
$ws.Range("F17").Value = 55188.31262192879
$ws.Range("G17").Value = 42123.38710481712
$ws.Range("H17").Value = 31060.40486387349
$ws.Range("B18").Value = 30757.77182006581
$ws.Range("C18").Value = 53117.0346995275
$ws.Range("D18").Value = 62588.97425345836
$ws.Range("E18").Value = 65712.52806844948
$ws.Range("F18").Value = 56177.35546690203
$ws.Range("G18").Value = 42859.54502267989
$ws.Range("H18").Value = 31421.04622032432
$ws.Range("B19").Value = 31372.84163022501
$ws.Range("C19").Value = 53836.09326920588
$ws.Range("D19").Value = 63501.33156767619
$ws.Range("E19").Value = 65937.37333986211
$ws.Range("F19").Value = 57299.70632113834
$ws.Range("G19").Value = 43612.83734516452
$ws.Range("H19").Value = 31814.2467355543
$ws.Range("B20").Value = 31770.65308312707
$ws.Range("C20").Value = 54396.9401941148
$ws.Range("D20").Value = 64332.00789307091
$ws.Range("E20").Value = 66388.84522073313
$ws.Range("F20").Value = 58488.7892248441
$ws.Range("G20").Value = 44444.36438733333
$ws.Range("H20").Value = 32281.25682786161
$ws.Range("B21").Value = 32085.87758099802
$ws.Range("C21").Value = 54895.78270432636
$ws.Range("D21").Value = 65172.58906450796
$ws.Range("E21").Value = 67074.61408725999
$ws.Range("F21").Value = 59581.67494743838
$ws.Range("G21").Value = 45457.14665086652
$ws.Range("H21").Value = 32833.83603827428
$ws.Range("B22").Value = 32604.43869893488
$ws.Range("C22").Value = 55424.93392210561
$ws.Range("D22").Value = 66252.90752221478
$ws.Range("E22").Value = 67893.17968626379
$ws.Range("F22").Value = 60610.96431462823
$ws.Range("G22").Value = 46597.77047270192
$ws.Range("H22").Value = 33442.45132945751
$ws.Range("B23").Value = 33203.23221482083
$ws.Range("C23").Value = 55861.14943500017
$ws.Range("D23").Value = 67342.34763425634
$ws.Range("E23").Value = 68599.29524819963
$ws.Range("F23").Value = 61417.3764847997
$ws.Range("G23").Value = 47637.29463226591
$ws.Range("H23").Value = 34114.21513137338
$ws.Range("B24").Value = 33579.73382334111
$ws.Range("C24").Value = 56003.11789673618
$ws.Range("D24").Value = 67967.31690799592
$ws.Range("E24").Value = 68942.38704144454
$ws.Range("F24").Value = 61879.86132062114
$ws.Range("G24").Value = 48393.37615892335
$ws.Range("H24").Value = 34810.14546411882
$ws.Range("B25").Value = 33652.57604712757
$ws.Range("C25").Value = 55712.53624805016
$ws.Range("D25").Value = 67888.70682097427
$ws.Range("E25").Value = 68848.02757543912
$ws.Range("F25").Value = 61955.70372876537
$ws.Range("G25").Value = 48758.8609105718
$ws.Range("H25").Value = 35422.31153244086
$ws.Range("B26").Value = 33432.19649475824
$ws.Range("C26").Value = 55111.14154910522
$ws.Range("D26").Value = 67139.33624375903
$ws.Range("E26").Value = 68368.40990551708
$ws.Range("F26").Value = 61547.16489863408
$ws.Range("G26").Value = 48822.8920107542
$ws.Range("H26").Value = 35866.26028534104
$ws.Range("B27").Value = 33102.25723697452
$ws.Range("C27").Value = 54368.09672727017
$ws.Range("D27").Value = 66014.31723196797
$ws.Range("E27").Value = 67547.88652745754
$ws.Range("F27").Value = 60793.41895337757
$ws.Range("G27").Value = 48850.01706455777
$ws.Range("H27").Value = 36125.21258736932
$ws.Range("B28").Value = 32918.06835922332
$ws.Range("C28").Value = 53737.2382652613
$ws.Range("D28").Value = 64935.09260604021
$ws.Range("E28").Value = 66800.39746065324
$ws.Range("F28").Value = 60113.15325806957
$ws.Range("G28").Value = 48960.68440505057
$ws.Range("H28").Value = 36310.70306265304
$ws.Range("B29").Value = 33004.57033905597
$ws.Range("C29").Value = 53314.82802815066
$ws.Range("D29").Value = 64274.95883254676
$ws.Range("E29").Value = 66513.21958885776
$ws.Range("F29").Value = 59831.13893948132
$ws.Range("G29").Value = 49361.39312198662
$ws.Range("H29").Value = 36583.82520623227
$ws.Range("B30").Value = 33208.98733407279
$ws.Range("C30").Value = 53068.79236474723
$ws.Range("D30").Value = 64172.85153146828
$ws.Range("E30").Value = 66883.28953732633
$ws.Range("F30").Value = 59894.78759116784
$ws.Range("G30").Value = 49950.59822863914
$ws.Range("H30").Value = 37028.08089905431
$ws.Range("B31").Value = 33422.6916507316
$ws.Range("C31").Value = 53001.65034871917
$ws.Range("D31").Value = 64690.02423348597
$ws.Range("E31").Value = 67809.476270536
$ws.Range("F31").Value = 60230.00788014476
$ws.Range("G31").Value = 50676.54071767023
$ws.Range("H31").Value = 37697.66451553714
$ws.Range("B32").Value = 33710.39644577842
$ws.Range("C32").Value = 53284.53029550634
$ws.Range("D32").Value = 65784.21092911209
$ws.Range("E32").Value = 69357.53238237416
$ws.Range("F32").Value = 60907.33114490833
$ws.Range("G32").Value = 51528.79851532031
$ws.Range("H32").Value = 38586.42320666646
$ws.Range("B33").Value = 34007.27536027006
$ws.Range("C33").Value = 53874.12459467207
$ws.Range("D33").Value = 67164.36496620174
$ws.Range("E33").Value = 71131.46628968956
$ws.Range("F33").Value = 61972.72096268774
$ws.Range("G33").Value = 52420.9749667956
$ws.Range("H33").Value = 39559.87189652432
$ws.Range("B34").Value = 34288.56680347348
$ws.Range("C34").Value = 54568.83864509473
$ws.Range("D34").Value = 68555.75554702959
$ws.Range("E34").Value = 72690.64283874938
$ws.Range("F34").Value = 63227.59396086118
$ws.Range("G34").Value = 53273.74432423768
$ws.Range("H34").Value = 40492.12280257665
$ws.Range("B35").Value = 34576.89272418692
$ws.Range("C35").Value = 55253.37823730725
$ws.Range("D35").Value = 69854.6472306981
$ws.Range("E35").Value = 73856.9672964909
$ws.Range("F35").Value = 64378.09114235226
$ws.Range("G35").Value = 53891.81696115074
$ws.Range("H35").Value = 41298.03578746622
$ws.Range("B36").Value = 35118.81323010771
$ws.Range("C36").Value = 55835.88704557788
$ws.Range("D36").Value = 70997.63392502052
$ws.Range("E36").Value = 74617.6551010385
$ws.Range("F36").Value = 65164.66749090133
$ws.Range("G36").Value = 53981.44480326099
$ws.Range("H36").Value = 41819.98899218474
$ws.Range("B37").Value = 36121.90767255844
$ws.Range("C37").Value = 56392.64307725892
$ws.Range("D37").Value = 71986.73006065389
$ws.Range("E37").Value = 75187.93764693913
$ws.Range("F37").Value = 65602.06020520019
$ws.Range("G37").Value = 53583.68763650253
$ws.Range("H37").Value = 42105.67653270228
$ws.Range("B38").Value = 37249.43326847029
$ws.Range("C38").Value = 57059.47650618546
$ws.Range("D38").Value = 72801.24748572899
$ws.Range("E38").Value = 75858.18165628714
$ws.Range("F38").Value = 65827.84381275695
$ws.Range("G38").Value = 53061.96578344086
$ws.Range("H38").Value = 42392.66577592886
$ws.Range("B39").Value = 38512.07756712458
$ws.Range("C39").Value = 57901.97536488777
$ws.Range("D39").Value = 73635.98152529351
$ws.Range("E39").Value = 77018.0967637063
$ws.Range("F39").Value = 66369.86217336563
$ws.Range("G39").Value = 52824.74855704819
$ws.Range("H39").Value = 43023.63547958494
$ws.Range("B40").Value = 39840.95286602573
$ws.Range("C40").Value = 58760.46709454913
$ws.Range("D40").Value = 74502.87080987265
$ws.Range("E40").Value = 78472.40841906633
$ws.Range("F40").Value = 67252.7916395301
$ws.Range("G40").Value = 52870.0075432985
$ws.Range("H40").Value = 43855.04043442155
$ws.Range("B41").Value = 41216.07560496016
$ws.Range("C41").Value = 59650.88898357334
$ws.Range("D41").Value = 75414.81985752781
$ws.Range("E41").Value = 80089.33059585026
$ws.Range("F41").Value = 68311.46739396745
$ws.Range("G41").Value = 53088.48166069209
$ws.Range("H41").Value = 44619.66404551171
